# TC10_CCDI_PHS-Accession-phs003111_StudyShrTitle-MolecCharClonal_LibStrat-WGS.xlsx
# Update the "startup" sheet's FilesTab row so the dbExcel/WebExcel input-file
# references point at the TC10 (this file's own) input excels instead of the
# TC09 ones that were previously referenced, and trim the FilesTab query's
# row limit from 100000 down to 100.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Activate()

$oldNeo4j = "TC09_CCDI_PHS-Accession-phs003111_Acro-CCDI_MSK_phs003111_LibSrc-Genomi_Neo4jData.xlsx"
$oldWeb   = "TC09_CCDI_PHS-Accession-phs003111_Acro-CCDI_MSK_phs003111_LibSrc-Genomi_WebData.xlsx"
$newNeo4j = "TC10_CCDI_PHS-Accession-phs003111_StudyShrTitle-MolecCharClonal_LibStrat-WGS_Neo4jData.xlsx"
$newWeb   = "TC10_CCDI_PHS-Accession-phs003111_StudyShrTitle-MolecCharClonal_LibStrat-WGS_WebData.xlsx"

# Column D = dbExcel (Neo4jData file name), Column E = WebExcel (WebData file name)
# Rows 2-6 correspond to ParticipantsTab, DiagnosisTab, StudiesTab, SamplesTab, FilesTab
for ($r = 2; $r -le 6; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    if ($dCell.Value() -eq $oldNeo4j) {
        $dCell.Value = $newNeo4j
    }
    if ($eCell.Value() -eq $oldWeb) {
        $eCell.Value = $newWeb
    }
}

# The FilesTab StatQuery (column C, row 6 / B6 in this sheet's layout) ends
# with "ORDER BY file_name LIMIT 100000" -- reduce the limit to 100.
$b6 = $ws.Range("B6")
$b6Text = $b6.Value()
$b6.Value = $b6Text.Replace("ORDER BY file_name LIMIT 100000", "ORDER BY file_name LIMIT 100")

# Move the saved selection from C18 to B12.
$ws.Range("B12").Select()
